$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "compression knee"
$ws.Range("A2").Value = "hockey pads men"
$ws.Range("A3").Value = "black basketball"
$ws.Range("A4").Value = "leg pads baseball"
$ws.Range("A5").Value = "running tights men"
$ws.Range("A6").Value = "snowboarding pads"
$ws.Range("A7").Value = "compression leggings youth boys"
$ws.Range("A8").Value = "knee pads wrestling men"
$ws.Range("A9").Value = "workout pants for girls"
$ws.Range("A10").Value = "basketball team gear"
$ws.Range("A11").Value = "basketball apparel for girls"
$ws.Range("A12").Value = "boys basketball compression leggings"
$ws.Range("A13").Value = "compressions pants"
$ws.Range("A14").Value = "bjj tights men"
$ws.Range("A15").Value = "compression basketball pants"
$ws.Range("A16").Value = "mens baseball pants"
$ws.Range("A17").Value = "best compression leggings"
$ws.Range("A18").Value = "softball pants black"
$ws.Range("A19").Value = "sliding pants youth"
$ws.Range("A20").Value = "basketballs pants"
$ws.Range("A21").Value = "sport knee pads"
$ws.Range("A22").Value = "men s knee pads"
$ws.Range("A23").Value = "volleyball knee pads black"
$ws.Range("A24").Value = "pack of leggings"
$ws.Range("A25").Value = "cycling knee protector"
$ws.Range("A26").Value = "knee pads sleeve basketball"
$ws.Range("A27").Value = "compression tights recovery"
$ws.Range("A28").Value = "volleyball pads"
$ws.Range("A29").Value = "padded tights for basketball"
$ws.Range("A30").Value = "small knee pads"
$ws.Range("A31").Value = "soccer apparel youth"
$ws.Range("A32").Value = "cheap knee pads"
$ws.Range("A33").Value = "baseball gear boys"
$ws.Range("A34").Value = "knee pads six six one"
$ws.Range("A35").Value = "compression pad"
$ws.Range("A36").Value = "free volleyball"
$ws.Range("A37").Value = "gym pad"
$ws.Range("A38").Value = "compression pants for boys basketball"
$ws.Range("A39").Value = "girl knee pads"
$ws.Range("A40").Value = "boys leggings sports youth"
$ws.Range("A41").Value = "youth football pants with pads small"
$ws.Range("A42").Value = "mens wrestling knee pads"
$ws.Range("A43").Value = "sport compression pants boys"
$ws.Range("A44").Value = "softball pants for girls"
$ws.Range("A45").Value = "compression leggings"
$ws.Range("A46").Value = "mens lacrosse pads"
$ws.Range("A47").Value = "softball pants youth girls black"
$ws.Range("A48").Value = "sport pants men"
$ws.Range("A49").Value = "thigh compression tights"
$ws.Range("A50").Value = "youth knee compression"
$ws.Range("A51").Value = "basketball sleeve knee pads"
$ws.Range("A52").Value = "leggings boys basketball"
$ws.Range("A53").Value = "basketball gear"
$ws.Range("A54").Value = "lacrosse youth pads"
$ws.Range("A55").Value = "running tights for men"
$ws.Range("A56").Value = "compressions leggings for men"
$ws.Range("A57").Value = "knee gel pads"
$ws.Range("A58").Value = "mens big and tall pants"
$ws.Range("A59").Value = "basketball protection"
$ws.Range("A60").Value = "knee braces for men xxl"
$ws.Range("A61").Value = "leggings with mesh girls"
$ws.Range("A62").Value = "little black pants guaranteed to fit"
$ws.Range("A63").Value = "3/4 pants men"
$ws.Range("A64").Value = "mens capri tights"
$ws.Range("A65").Value = "boys basketball"
$ws.Range("A66").Value = "girls volleyball spandex"
$ws.Range("A67").Value = "girls softball pants"
$ws.Range("A68").Value = "lcl knee support"
$ws.Range("A69").Value = "arthritis equipment"
$ws.Range("A70").Value = "volleyball knee pads"
$ws.Range("A71").Value = "running compression men"
$ws.Range("A72").Value = "baseball pants youth xl"
$ws.Range("A73").Value = "youth xs knee pads basketball"
$ws.Range("A74").Value = "under pants for basketball"
$ws.Range("A75").Value = "guy legging"
$ws.Range("A76").Value = "under armor warm pants"
$ws.Range("A77").Value = "spandex leggings men thermal"
$ws.Range("A78").Value = "mems thermal leggings"
$ws.Range("A79").Value = "men compression 3 4 pants"
$ws.Range("A80").Value = "men compression leggings 3 4"
$ws.Range("A81").Value = "men compression pants adidas"
$ws.Range("A82").Value = "men compression tights 3 4"
$ws.Range("A83").Value = "mens compression 3 4 tights"
$ws.Range("A84").Value = "mens compression 3 4 leggings"
$ws.Range("A85").Value = "eastbay baseball pants"
$ws.Range("A86").Value = "eastbay leggings"
$ws.Range("A87").Value = "eastbay tights"
$ws.Range("A88").Value = "mcdavid knee pads youth"
$ws.Range("A89").Value = "basketball leggings nike"
$ws.Range("A90").Value = "nba basketball pants"
$ws.Range("A91").Value = "basketball 3 4 leggings"
$ws.Range("A92").Value = "elbow knee wrist pads for youth"
$ws.Range("A93").Value = "telsa thermals men"
$ws.Range("A94").Value = "wintergear compression leggings men"
$ws.Range("A95").Value = "dry skin tights"
$ws.Range("A96").Value = "mens workoit tights"
$ws.Range("A97").Value = "track leggings"
$ws.Range("A98").Value = "track tights for men"
$ws.Range("A99").Value = "tesla wintergear compression shirt"
$ws.Range("A100").Value = "basketball apparel youth"
